$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows right after row 25 (PERIOD row for Feb 2023).
# This pushes the existing monthly PERIOD rows (and everything else below)
# down by two rows, matching the diff where row26's old date (3/1/2023)
# ends up on row28, etc.
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

# The freshly inserted rows 26-27 don't carry the normal data-row
# formatting, so copy it over from row 28 (which still has the original
# "blank period row" formatting/style indexes).
$ws.Range("A28:K28").Copy()
$ws.Range("A26:K27").PasteSpecial(-4122)

$earnedFormula = 'IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G26").Formula = "=" + $earnedFormula
$ws.Range("G27").Formula = "=" + $earnedFormula

# Populate the new leave entries.
$ws.Range("B27").Value = "VL(12-0-0)"
$ws.Range("F27").Value = 12
$ws.Range("K27").Value = "CAMBODIA 5/3 - 18"
$ws.Range("K26").Value = "BALI INDONESIA 4/29 - 5/1"

# Expand the table so it covers the two new rows added at the bottom
# (Excel keeps the table the same size unless we explicitly grow it),
# keeping the special "last row" formatting on the new final row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K137"))
$ws.Range("G136").Formula = "=" + $earnedFormula
$ws.Range("G137").Formula = "=" + $earnedFormula

# Update the active selection to match the saved state of the workbook.
$ws.Range("B26").Select()
